$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3506.2632  # H64 (was 3507.4805)
$ws.Cells.Item(64, 9).Value = 3594.4211  # I64 (was 3604.9473)
$ws.Cells.Item(64, 10).Value = 3241.7896  # J64 (was 3229.7)
$ws.Cells.Item(64, 11).Value = 3594.4211  # K64 (was 3604.9473)
$ws.Cells.Item(64, 12).Value = 3241.7896  # L64 (was 3229.7)
$ws.Cells.Item(64, 13).Value = -3346.4211  # M64 (was -3356.9473)
$ws.Cells.Item(64, 14).Value = -3737.7896  # N64 (was -3725.7)
$ws.Cells.Item(67, 8).Value = 3506.2632  # H67 (was 3507.4805)
$ws.Cells.Item(67, 9).Value = 3594.4211  # I67 (was 3604.9473)
$ws.Cells.Item(67, 10).Value = 3241.7896  # J67 (was 3229.7)
$ws.Cells.Item(67, 11).Value = 3594.4211  # K67 (was 3604.9473)
$ws.Cells.Item(67, 12).Value = 3241.7896  # L67 (was 3229.7)
$ws.Cells.Item(67, 13).Value = -2736.4211  # M67 (was -2746.9473)
$ws.Cells.Item(67, 14).Value = -4957.7896  # N67 (was -4945.7)
$ws.Cells.Item(82, 8).Value = 1165.3334  # H82 (was 1049.2)
$ws.Cells.Item(82, 9).Value = 1165.3334  # I82 (was 1348.6666)
$ws.Cells.Item(82, 10).Value = 0  # J82 (was 600)
$ws.Cells.Item(82, 11).Value = 3496.0002  # K82 (was 4045.9998)
$ws.Cells.Item(82, 12).Value = 0  # L82 (was 1800)
$ws.Cells.Item(82, 13).Value = -3090.0002  # M82 (was -3639.9998)
$ws.Cells.Item(82, 14).Value = $null  # N82 delete (was -2612)
$ws.Cells.Item(85, 8).Value = 1165.3334  # H85 (was 1049.2)
$ws.Cells.Item(85, 9).Value = 1165.3334  # I85 (was 1348.6666)
$ws.Cells.Item(85, 10).Value = 0  # J85 (was 600)
$ws.Cells.Item(85, 11).Value = 3496.0002  # K85 (was 4045.9998)
$ws.Cells.Item(85, 12).Value = 0  # L85 (was 1800)
$ws.Cells.Item(85, 13).Value = -2092.0002  # M85 (was -2641.9998)
$ws.Cells.Item(85, 14).Value = $null  # N85 delete (was -4608)
$ws.Cells.Item(86, 8).Value = 2224.4666  # H86 (was 2221.6875)
$ws.Cells.Item(86, 9).Value = 4901  # I86 (was 2600.5)
$ws.Cells.Item(86, 10).Value = 1555.3334  # J86 (was 1994.4)
$ws.Cells.Item(86, 11).Value = 4901  # K86 (was 2600.5)
$ws.Cells.Item(86, 12).Value = 1555.3334  # L86 (was 1994.4)
$ws.Cells.Item(86, 13).Value = -3778  # M86 (was -1477.5)
$ws.Cells.Item(86, 14).Value = -3801.3334  # N86 (was -4240.4)
$ws.Cells.Item(89, 8).Value = 2224.4666  # H89 (was 2221.6875)
$ws.Cells.Item(89, 9).Value = 4901  # I89 (was 2600.5)
$ws.Cells.Item(89, 10).Value = 1555.3334  # J89 (was 1994.4)
$ws.Cells.Item(89, 11).Value = 24505  # K89 (was 13002.5)
$ws.Cells.Item(89, 12).Value = 7776.666999999999  # L89 (was 9972)
$ws.Cells.Item(89, 13).Value = -18889  # M89 (was -7386.5)
$ws.Cells.Item(89, 14).Value = -19008.667  # N89 (was -21204)
$ws.Cells.Item(92, 8).Value = 49604016  # H92 (was 49604060)
$ws.Cells.Item(92, 9).Value = 2778576.2  # I92 (was 2924860)
$ws.Cells.Item(92, 10).Value = 166667620  # J92 (was 148149040)
$ws.Cells.Item(92, 11).Value = 2778576.2  # K92 (was 2924860)
$ws.Cells.Item(92, 12).Value = 166667620  # L92 (was 148149040)
$ws.Cells.Item(92, 13).Value = -2777328.2  # M92 (was -2923612)
$ws.Cells.Item(92, 14).Value = -166670116  # N92 (was -148151536)
$ws.Cells.Item(112, 8).Value = 4202690.5  # H112 (was 3969285.5)
$ws.Cells.Item(112, 10).Value = 5103102.5  # J112 (was 4762989)
$ws.Cells.Item(112, 12).Value = 15309307.5  # L112 (was 14288967)
$ws.Cells.Item(112, 14).Value = -15311523.5  # N112 (was -14291183)
$ws.Cells.Item(132, 8).Value = 1111.5094  # H132 (was 1168.7693)
$ws.Cells.Item(132, 9).Value = 980.08  # I132 (was 1022.1875)
$ws.Cells.Item(132, 10).Value = 3302  # J132 (was 2927.75)
$ws.Cells.Item(132, 11).Value = 2940.24  # K132 (was 3066.5625)
$ws.Cells.Item(132, 12).Value = 9906  # L132 (was 8783.25)
$ws.Cells.Item(132, 13).Value = -410.2400000000002  # M132 (was -536.5625)
$ws.Cells.Item(132, 14).Value = -14966  # N132 (was -13843.25)
$ws.Cells.Item(138, 8).Value = 1460.238  # H138 (was 1611.921)
$ws.Cells.Item(138, 9).Value = 1147.7084  # I138 (was 1332.4)
$ws.Cells.Item(138, 10).Value = 1876.9445  # J138 (was 1922.5)
$ws.Cells.Item(138, 11).Value = 3443.1252  # K138 (was 3997.2)
$ws.Cells.Item(138, 12).Value = 5630.833500000001  # L138 (was 5767.5)
$ws.Cells.Item(138, 13).Value = 1696.8748  # M138 (was 1142.8)
$ws.Cells.Item(138, 14).Value = -15910.8335  # N138 (was -16047.5)

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 8813.294  # H61 (was 6029.44)
$ws.Cells.Item(61, 9).Value = 10368.833  # I61 (was 6368.8)
$ws.Cells.Item(61, 10).Value = 5080  # J61 (was 4672)
$ws.Cells.Item(61, 11).Value = 10368.833  # K61 (was 6368.8)
$ws.Cells.Item(61, 12).Value = 5080  # L61 (was 4672)
$ws.Cells.Item(61, 13).Value = -10156.833  # M61 (was -6156.8)
$ws.Cells.Item(61, 14).Value = -5504  # N61 (was -5096)
$ws.Cells.Item(74, 8).Value = 2112.6924  # H74 (was 2213.75)
$ws.Cells.Item(74, 9).Value = 1935.8889  # I74 (was 2065.375)
$ws.Cells.Item(74, 11).Value = 1935.8889  # K74 (was 2065.375)
$ws.Cells.Item(74, 13).Value = -1061.8889  # M74 (was -1191.375)
$ws.Cells.Item(77, 8).Value = 2112.6924  # H77 (was 2213.75)
$ws.Cells.Item(77, 9).Value = 1935.8889  # I77 (was 2065.375)
$ws.Cells.Item(77, 11).Value = 9679.4445  # K77 (was 10326.875)
$ws.Cells.Item(77, 13).Value = -5311.4445  # M77 (was -5958.875)
$ws.Cells.Item(102, 8).Value = 3705593.5  # H102 (was 5292976.5)
$ws.Cells.Item(102, 9).Value = 4116992.8  # I102 (was 6174639)
$ws.Cells.Item(102, 11).Value = 4116992.8  # K102 (was 6174639)
$ws.Cells.Item(102, 13).Value = -4115370.8  # M102 (was -6173017)
$ws.Cells.Item(132, 8).Value = 2195.5615  # H132 (was 2260.6)
$ws.Cells.Item(132, 9).Value = 1354.6818  # I132 (was 1399.8096)
$ws.Cells.Item(132, 11).Value = 4064.0454  # K132 (was 4199.4288)
$ws.Cells.Item(132, 13).Value = -1534.0454  # M132 (was -1669.4288)
$ws.Cells.Item(136, 8).Value = 8813.294  # H136 (was 6029.44)
$ws.Cells.Item(136, 9).Value = 10368.833  # I136 (was 6368.8)
$ws.Cells.Item(136, 10).Value = 5080  # J136 (was 4672)
$ws.Cells.Item(136, 11).Value = 31106.499  # K136 (was 19106.4)
$ws.Cells.Item(136, 12).Value = 15240  # L136 (was 14016)
$ws.Cells.Item(136, 13).Value = -28556.499  # M136 (was -16556.4)
$ws.Cells.Item(136, 14).Value = -20340  # N136 (was -19116)

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1836.909  # H86 (was 1883.8334)
$ws.Cells.Item(86, 9).Value = 1800.6666  # I86 (was 1845.1111)
$ws.Cells.Item(86, 11).Value = 1800.6666  # K86 (was 1845.1111)
$ws.Cells.Item(86, 13).Value = -677.6666  # M86 (was -722.1111000000001)
$ws.Cells.Item(89, 8).Value = 1836.909  # H89 (was 1883.8334)
$ws.Cells.Item(89, 9).Value = 1800.6666  # I89 (was 1845.1111)
$ws.Cells.Item(89, 11).Value = 9003.333000000001  # K89 (was 9225.5555)
$ws.Cells.Item(89, 13).Value = -3387.333000000001  # M89 (was -3609.5555)
$ws.Cells.Item(94, 8).Value = 1858.48  # H94 (was 1906.875)
$ws.Cells.Item(94, 9).Value = 1400.2858  # I94 (was 1414.7858)
$ws.Cells.Item(94, 10).Value = 2441.6365  # J94 (was 2595.8)
$ws.Cells.Item(94, 11).Value = 1400.2858  # K94 (was 1414.7858)
$ws.Cells.Item(94, 12).Value = 2441.6365  # L94 (was 2595.8)
$ws.Cells.Item(94, 13).Value = -949.2858000000001  # M94 (was -963.7858000000001)
$ws.Cells.Item(94, 14).Value = -3343.6365  # N94 (was -3497.8)
$ws.Cells.Item(99, 8).Value = 38462744  # H99 (was 40001240)
$ws.Cells.Item(99, 9).Value = 52632696  # I99 (was 55556710)
$ws.Cells.Item(99, 11).Value = 52632696  # K99 (was 55556710)
$ws.Cells.Item(99, 13).Value = -52631198  # M99 (was -55555212)
$ws.Cells.Item(134, 8).Value = 4320.857  # H134 (was 4247.5815)
$ws.Cells.Item(134, 9).Value = 4936.3228  # I134 (was 4803)
$ws.Cells.Item(134, 10).Value = 2586.3635  # J134 (was 2631.818)
$ws.Cells.Item(134, 11).Value = 14808.9684  # K134 (was 14409)
$ws.Cells.Item(134, 12).Value = 7759.0905  # L134 (was 7895.454000000001)
$ws.Cells.Item(134, 13).Value = -12273.9684  # M134 (was -11874)
$ws.Cells.Item(134, 14).Value = -12829.0905  # N134 (was -12965.454)

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1249  # H132 (was 1350)
$ws.Cells.Item(132, 9).Value = 1008.8571  # I132 (was 1100.9778)
$ws.Cells.Item(132, 10).Value = 2318.7273  # J132 (was 2470.6)
$ws.Cells.Item(132, 11).Value = 3026.5713  # K132 (was 3302.9334)
$ws.Cells.Item(132, 12).Value = 6956.1819  # L132 (was 7411.799999999999)
$ws.Cells.Item(132, 13).Value = -496.5712999999996  # M132 (was -772.9333999999999)
$ws.Cells.Item(132, 14).Value = -12016.1819  # N132 (was -12471.8)
$ws.Cells.Item(134, 8).Value = 1909.4482  # H134 (was 1964.2142)
$ws.Cells.Item(134, 9).Value = 2074.4783  # I134 (was 2151.682)
$ws.Cells.Item(134, 11).Value = 6223.4349  # K134 (was 6455.045999999999)
$ws.Cells.Item(134, 13).Value = -3688.4349  # M134 (was -3920.045999999999)

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2128550.5  # H131 (was 2174802.5)
$ws.Cells.Item(131, 10).Value = 986.8461  # J131 (was 987.5526)
$ws.Cells.Item(131, 12).Value = 2960.5383  # L131 (was 2962.6578)
$ws.Cells.Item(131, 14).Value = -13040.5383  # N131 (was -13042.6578)

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2774.6191  # H132 (was 2856.6667)
$ws.Cells.Item(132, 9).Value = 2252.7  # I132 (was 2690.2856)
$ws.Cells.Item(132, 10).Value = 3249.0908  # J132 (was 2939.8572)
$ws.Cells.Item(132, 11).Value = 6758.099999999999  # K132 (was 8070.8568)
$ws.Cells.Item(132, 12).Value = 9747.2724  # L132 (was 8819.571599999999)
$ws.Cells.Item(132, 13).Value = -4228.099999999999  # M132 (was -5540.8568)
$ws.Cells.Item(132, 14).Value = -14807.2724  # N132 (was -13879.5716)

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1207.25  # H93 (was 1207.5)
$ws.Cells.Item(93, 10).Value = 1276.7778  # J93 (was 1277.1111)
$ws.Cells.Item(93, 12).Value = 1276.7778  # L93 (was 1277.1111)
$ws.Cells.Item(93, 14).Value = -3772.7778  # N93 (was -3773.1111)
$ws.Cells.Item(132, 8).Value = 12826274  # H132 (was 11910104)
$ws.Cells.Item(132, 9).Value = 18525062  # I132 (was 16672546)
$ws.Cells.Item(132, 11).Value = 55575186  # K132 (was 50017638)
$ws.Cells.Item(132, 13).Value = -55572656  # M132 (was -50015108)
$ws.Cells.Item(136, 8).Value = 7026.5835  # H136 (was 6896.1353)
$ws.Cells.Item(136, 9).Value = 5681.7334  # I136 (was 5690.067)
$ws.Cells.Item(136, 10).Value = 13750.833  # J136 (was 12065)
$ws.Cells.Item(136, 11).Value = 17045.2002  # K136 (was 17070.201)
$ws.Cells.Item(136, 12).Value = 41252.499  # L136 (was 36195)
$ws.Cells.Item(136, 13).Value = -14495.2002  # M136 (was -14520.201)
$ws.Cells.Item(136, 14).Value = -46352.499  # N136 (was -41295)

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1322.8959  # H132 (was 1312.3877)
$ws.Cells.Item(132, 9).Value = 984.43335  # I132 (was 984.7)
$ws.Cells.Item(132, 10).Value = 1887  # J132 (was 1829.7894)
$ws.Cells.Item(132, 11).Value = 2953.30005  # K132 (was 2954.1)
$ws.Cells.Item(132, 12).Value = 5661  # L132 (was 5489.3682)
$ws.Cells.Item(132, 13).Value = -423.3000499999998  # M132 (was -424.1000000000004)
$ws.Cells.Item(132, 14).Value = -10721  # N132 (was -10549.3682)
$ws.Cells.Item(136, 8).Value = 1925.6316  # H136 (was 2049.6)
$ws.Cells.Item(136, 9).Value = 1845.7878  # I136 (was 2175.3215)
$ws.Cells.Item(136, 10).Value = 2035.4166  # J136 (was 1919.2222)
$ws.Cells.Item(136, 11).Value = 5537.3634  # K136 (was 6525.9645)
$ws.Cells.Item(136, 12).Value = 6106.2498  # L136 (was 5757.6666)
$ws.Cells.Item(136, 13).Value = -2987.3634  # M136 (was -3975.9645)
$ws.Cells.Item(136, 14).Value = -11206.2498  # N136 (was -10857.6666)

